$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A=3; B="Femacal de La Calera"; C="Coquimbo"; D=45008; E=5; F="Fruta"; G=100104; H="Frutos de pepita"; I=100104003; J="Membrillo"; K="Champion"; L="Especial";                 M=67; N=18000; O=18000; P=18000; Q="`$/caja 18 kilos empedrada"; R="Región de O'Higgins"; S=1000; T=18 },
    @{ A=3; B="Femacal de La Calera"; C="Coquimbo"; D=45008; E=5; F="Fruta"; G=100104; H="Frutos de pepita"; I=100104003; J="Membrillo"; K="Champion"; L="Extra (doble especial)";   M=56; N=20000; O=20000; P=20000; Q="`$/caja 18 kilos empedrada"; R="Región de O'Higgins"; S=1111; T=18 },
    @{ A=3; B="Femacal de La Calera"; C="Coquimbo"; D=45008; E=5; F="Fruta"; G=100104; H="Frutos de pepita"; I=100104003; J="Membrillo"; K="Champion"; L="Primera";                   M=60; N=16000; O=16000; P=16000; Q="`$/caja 18 kilos empedrada"; R="Región de O'Higgins"; S=889;  T=18 }
)

$startRow = 100
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
    $ws.Cells.Item($r, 11).Value = $data.K
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 17).Value = $data.Q
    $ws.Cells.Item($r, 18).Value = $data.R
    $ws.Cells.Item($r, 19).Value = $data.S
    $ws.Cells.Item($r, 20).Value = $data.T
}
